$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3731191
$ws.Range("J17").Value = 3731191
$ws.Range("L17").Value = 11193573
$ws.Range("N17").Value = -11193909

$ws.Range("H87").Value = 30899.5
$ws.Range("J87").Value = 30899.5
$ws.Range("L87").Value = 30899.5
$ws.Range("N87").Value = -33395.5

$ws.Range("H90").Value = 30899.5
$ws.Range("J90").Value = 30899.5
$ws.Range("L90").Value = 92698.5
$ws.Range("N90").Value = -105178.5

$ws.Range("H127").Value = 1667.6976
$ws.Range("I127").Value = 351.75
$ws.Range("J127").Value = 1968.4857
$ws.Range("K127").Value = 1055.25
$ws.Range("L127").Value = 5905.4571
$ws.Range("M127").Value = 3904.75
$ws.Range("N127").Value = -15825.4571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 42300
$ws.Range("J62").Value = 42300
$ws.Range("L62").Value = 42300
$ws.Range("N62").Value = -43548

$ws.Range("H64").Value = 42550
$ws.Range("J64").Value = 42550
$ws.Range("L64").Value = 42550
$ws.Range("N64").Value = -43046

$ws.Range("H65").Value = 42300
$ws.Range("J65").Value = 42300
$ws.Range("L65").Value = 126900
$ws.Range("N65").Value = -133140

$ws.Range("H67").Value = 42550
$ws.Range("J67").Value = 42550
$ws.Range("L67").Value = 42550
$ws.Range("N67").Value = -44266

$ws.Range("H75").Value = 40086.5
$ws.Range("J75").Value = 40086.5
$ws.Range("L75").Value = 40086.5
$ws.Range("N75").Value = -41834.5

$ws.Range("H76").Value = 31285.715
$ws.Range("J76").Value = 31285.715
$ws.Range("L76").Value = 31285.715
$ws.Range("N76").Value = -31961.715

$ws.Range("H78").Value = 40086.5
$ws.Range("J78").Value = 40086.5
$ws.Range("L78").Value = 120259.5
$ws.Range("N78").Value = -128995.5

$ws.Range("H79").Value = 31285.715
$ws.Range("J79").Value = 31285.715
$ws.Range("L79").Value = 31285.715
$ws.Range("N79").Value = -33625.715

$ws.Range("H81").Value = 42200
$ws.Range("J81").Value = 42200
$ws.Range("L81").Value = 42200
$ws.Range("N81").Value = -44196

$ws.Range("H84").Value = 42200
$ws.Range("J84").Value = 42200
$ws.Range("L84").Value = 126600
$ws.Range("N84").Value = -136584

$ws.Range("H103").Value = 49444
$ws.Range("J103").Value = 49444
$ws.Range("L103").Value = 49444
$ws.Range("N103").Value = -51788

$ws.Range("H132").Value = 2516.5
$ws.Range("I132").Value = 1525.0883
$ws.Range("J132").Value = 5325.5
$ws.Range("K132").Value = 4575.2649
$ws.Range("L132").Value = 15976.5
$ws.Range("M132").Value = -2045.2649
$ws.Range("N132").Value = -21036.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 42197.5
$ws.Range("J68").Value = 42197.5
$ws.Range("L68").Value = 42197.5
$ws.Range("N68").Value = -43819.5

$ws.Range("H71").Value = 42197.5
$ws.Range("J71").Value = 42197.5
$ws.Range("L71").Value = 126592.5
$ws.Range("N71").Value = -134704.5

$ws.Range("H75").Value = 24200
$ws.Range("I75").Value = 1000
$ws.Range("K75").Value = 1000
$ws.Range("M75").Value = -64

$ws.Range("H78").Value = 24200
$ws.Range("I78").Value = 1000
$ws.Range("K78").Value = 3000
$ws.Range("M78").Value = 1680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3919.3872
$ws.Range("I99").Value = 3000.0344
$ws.Range("K99").Value = 3000.0344
$ws.Range("M99").Value = -1502.0344

$ws.Range("H126").Value = 3919.3872
$ws.Range("I126").Value = 3000.0344
$ws.Range("K126").Value = 9000.1032
$ws.Range("M126").Value = -6530.1032

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1333
$ws.Range("I34").Value = 450
$ws.Range("J34").Value = 1518.8948
$ws.Range("K34").Value = 1350
$ws.Range("L34").Value = 4556.6844
$ws.Range("M34").Value = -1266
$ws.Range("N34").Value = -4724.6844

$ws.Range("H39").Value = 1125
$ws.Range("J39").Value = 1125
$ws.Range("L39").Value = 3375
$ws.Range("N39").Value = -3963

$ws.Range("H55").Value = 2259.9333
$ws.Range("J55").Value = 2259.9333
$ws.Range("L55").Value = 6779.7999
$ws.Range("N55").Value = -7133.7999

$ws.Range("H97").Value = 16667280
$ws.Range("I97").Value = 33333666
$ws.Range("J97").Value = 893.3333
$ws.Range("K97").Value = 100000998
$ws.Range("L97").Value = 2679.9999
$ws.Range("M97").Value = -100000502
$ws.Range("N97").Value = -3671.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5790.909
$ws.Range("J33").Value = 5790.909
$ws.Range("L33").Value = 5790.909
$ws.Range("N33").Value = -6294.909

$ws.Range("H52").Value = 24750
$ws.Range("J52").Value = 24750
$ws.Range("L52").Value = 24750
$ws.Range("N52").Value = -25268

$ws.Range("H64").Value = 32282.445
$ws.Range("J64").Value = 32282.445
$ws.Range("L64").Value = 32282.445
$ws.Range("N64").Value = -32778.445

$ws.Range("H67").Value = 32282.445
$ws.Range("J67").Value = 32282.445
$ws.Range("L67").Value = 32282.445
$ws.Range("N67").Value = -33998.445

$ws.Range("H69").Value = 38150.75
$ws.Range("J69").Value = 38150.75
$ws.Range("L69").Value = 38150.75
$ws.Range("N69").Value = -39648.75

$ws.Range("H72").Value = 38150.75
$ws.Range("J72").Value = 38150.75
$ws.Range("L72").Value = 114452.25
$ws.Range("N72").Value = -121940.25

$ws.Range("H82").Value = 32000
$ws.Range("J82").Value = 32000
$ws.Range("L82").Value = 32000
$ws.Range("N82").Value = -32766

$ws.Range("H85").Value = 32000
$ws.Range("J85").Value = 32000
$ws.Range("L85").Value = 32000
$ws.Range("N85").Value = -34652

$ws.Range("H132").Value = 8336096
$ws.Range("I132").Value = 11113461
$ws.Range("K132").Value = 33340383
$ws.Range("M132").Value = -33337853

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1731.4286
$ws.Range("I7").Value = 1586.6666
$ws.Range("J7").Value = 2600
$ws.Range("K7").Value = 1586.6666
$ws.Range("L7").Value = 2600
$ws.Range("M7").Value = -1474.6666
$ws.Range("N7").Value = -2824

$ws.Range("H46").Value = 10101645
$ws.Range("I46").Value = 66667096
$ws.Range("J46").Value = 672.0714
$ws.Range("K46").Value = 66667096
$ws.Range("L46").Value = 672.0714
$ws.Range("M46").Value = -66666908
$ws.Range("N46").Value = -1048.0714

$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25450

$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26560

$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H70").Value = 11500
$ws.Range("J70").Value = 11500
$ws.Range("L70").Value = 11500
$ws.Range("N70").Value = -12040

$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H73").Value = 11500
$ws.Range("J73").Value = 11500
$ws.Range("L73").Value = 11500
$ws.Range("N73").Value = -13372

$ws.Range("H126").Value = 1731.4286
$ws.Range("I126").Value = 1586.6666
$ws.Range("J126").Value = 2600
$ws.Range("K126").Value = 4759.9998
$ws.Range("L126").Value = 7800
$ws.Range("M126").Value = -2289.9998
$ws.Range("N126").Value = -12740

$ws.Range("H132").Value = 15811149
$ws.Range("I132").Value = 20840660
$ws.Range("J132").Value = 4112.143
$ws.Range("K132").Value = 62521980
$ws.Range("L132").Value = 12336.429
$ws.Range("M132").Value = -62519450
$ws.Range("N132").Value = -17396.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2518.5881
$ws.Range("I132").Value = 2001.8572
$ws.Range("J132").Value = 2880.3
$ws.Range("K132").Value = 6005.571599999999
$ws.Range("L132").Value = 8640.900000000001
$ws.Range("M132").Value = -3475.571599999999
$ws.Range("N132").Value = -13700.9
